$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gip"
$ws.Cells.Item(2, 3).Value = "Gipr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3260275
$ws.Cells.Item(2, 8).Value = 0.6520550000000001
$ws.Cells.Item(2, 9).Value = 0.4722113996121241
$ws.Cells.Item(2, 10).Value = 0.4126724043544658
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.205225
$ws.Cells.Item(2, 14).Value = 0.41045
$ws.Cells.Item(2, 15).Value = 0.3926590127895196
$ws.Cells.Item(2, 16).Value = 0.3926590127895196
$ws.Cells.Item(2, 17).Value = 0.0669089936875
$ws.Cells.Item(2, 18).Value = 0.26763597475
$ws.Cells.Item(2, 19).Value = 0.185418061999654
$ws.Cells.Item(2, 20).Value = 0.162039538899302

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gip"
$ws.Cells.Item(3, 3).Value = "Gipr"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3260275
$ws.Cells.Item(3, 8).Value = 0.6520550000000001
$ws.Cells.Item(3, 9).Value = 0.4722113996121241
$ws.Cells.Item(3, 10).Value = 0.4126724043544658
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3174295
$ws.Cells.Item(3, 14).Value = 0.6348590000000001
$ws.Cells.Item(3, 15).Value = 0.6073409872104804
$ws.Cells.Item(3, 16).Value = 0.6073409872104804
$ws.Cells.Item(3, 17).Value = 0.10349074631125
$ws.Cells.Item(3, 18).Value = 0.4139629852450001
$ws.Cells.Item(3, 19).Value = 0.2867933376124701
$ws.Cells.Item(3, 20).Value = 0.2506328654551638

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gip"
$ws.Cells.Item(4, 3).Value = "Gipr"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.199225
$ws.Cells.Item(4, 8).Value = 0.597675
$ws.Cells.Item(4, 9).Value = 0.2885533155568945
$ws.Cells.Item(4, 10).Value = 0.3782564036355144
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.5
$ws.Cells.Item(4, 13).Value = 0.205225
$ws.Cells.Item(4, 14).Value = 0.41045
$ws.Cells.Item(4, 15).Value = 0.3926590127895196
$ws.Cells.Item(4, 16).Value = 0.3926590127895196
$ws.Cells.Item(4, 17).Value = 0.04088595062499999
$ws.Cells.Item(4, 18).Value = 0.24531570375
$ws.Cells.Item(4, 19).Value = 0.1133030600237129
$ws.Cells.Item(4, 20).Value = 0.1485257860328351

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gip"
$ws.Cells.Item(5, 3).Value = "Gipr"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.199225
$ws.Cells.Item(5, 8).Value = 0.597675
$ws.Cells.Item(5, 9).Value = 0.2885533155568945
$ws.Cells.Item(5, 10).Value = 0.3782564036355144
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3174295
$ws.Cells.Item(5, 14).Value = 0.6348590000000001
$ws.Cells.Item(5, 15).Value = 0.6073409872104804
$ws.Cells.Item(5, 16).Value = 0.6073409872104804
$ws.Cells.Item(5, 17).Value = 0.0632398921375
$ws.Cells.Item(5, 18).Value = 0.379439352825
$ws.Cells.Item(5, 19).Value = 0.1752502555331815
$ws.Cells.Item(5, 20).Value = 0.2297306176026792

$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Gip"
$ws.Cells.Item(6, 3).Value = "Gipr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.5
$ws.Cells.Item(6, 7).Value = 0.1651745
$ws.Cells.Item(6, 8).Value = 0.330349
$ws.Cells.Item(6, 9).Value = 0.2392352848309814
$ws.Cells.Item(6, 10).Value = 0.2090711920100198
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.5
$ws.Cells.Item(6, 13).Value = 0.205225
$ws.Cells.Item(6, 14).Value = 0.41045
$ws.Cells.Item(6, 15).Value = 0.3926590127895196
$ws.Cells.Item(6, 16).Value = 0.3926590127895196
$ws.Cells.Item(6, 17).Value = 0.0338979367625
$ws.Cells.Item(6, 18).Value = 0.13559174705
$ws.Cells.Item(6, 19).Value = 0.09393789076615271
$ws.Cells.Item(6, 20).Value = 0.08209368785738247

$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Gip"
$ws.Cells.Item(7, 3).Value = "Gipr"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.5
$ws.Cells.Item(7, 7).Value = 0.1651745
$ws.Cells.Item(7, 8).Value = 0.330349
$ws.Cells.Item(7, 9).Value = 0.2392352848309814
$ws.Cells.Item(7, 10).Value = 0.2090711920100198
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3174295
$ws.Cells.Item(7, 14).Value = 0.6348590000000001
$ws.Cells.Item(7, 15).Value = 0.6073409872104804
$ws.Cells.Item(7, 16).Value = 0.6073409872104804
$ws.Cells.Item(7, 17).Value = 0.05243125894775
$ws.Cells.Item(7, 18).Value = 0.209725035791
$ws.Cells.Item(7, 19).Value = 0.1452973940648287
$ws.Cells.Item(7, 20).Value = 0.1269775041526373

$ws.Rows("8:10").Delete()
Write-Output "OK"
